# EIA Table 4.18 monthly refresh: October 2016 -> November 2016 data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title text update (A2, merged A2:J2)
$ws.Range("A2").Value = "Industrial Sector by State, November 2016"

# --- New England ---
$ws.Range("B5").Value = 2          # Row 5  (New England)
$ws.Range("B7").Value = 2          # Row 7  (Maine)

# --- Middle Atlantic ---
$ws.Range("B12").Value = 25        # Row 12 (Middle Atlantic)
$ws.Range("C12").Value = 1.71
$ws.Range("D12").Value = 8.6

$ws.Range("B14").Value = 24        # Row 14 (New York)
$ws.Range("C14").Value = 1.71
$ws.Range("D14").Value = 8.6

$ws.Range("B15").Value = 0.08      # Row 15 (Pennsylvania)
$ws.Range("C15").Value = 1.6
$ws.Range("D15").Value = 11.7

# --- East North Central ---
$ws.Range("B16").Value = 111       # Row 16 (East North Central)
$ws.Range("C16").Value = 3.3
$ws.Range("D16").Value = 8.9
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 0.62

$ws.Range("B17").Value = 98        # Row 17 (Illinois)
$ws.Range("E17").Value = 49
$ws.Range("F17").Value = 0.8
$ws.Range("G17").Value = 6.5

# --- West North Central ---
$ws.Range("D21").Value = 11.5      # Row 21 (Missouri)
$ws.Range("E21").Value = 26
$ws.Range("F21").Value = 0.28000000000000003
$ws.Range("G21").Value = 5.8

$ws.Range("B22").Value = 11        # Row 22 (South Atlantic)
$ws.Range("D22").Value = 7.4
$ws.Range("E22").Value = 49

$ws.Range("B23").Value = 11        # Row 23 (Delaware)
$ws.Range("D23").Value = 7.4
$ws.Range("E23").Value = 49

# --- East South Central ---
$ws.Range("B30").Value = 57        # Row 30 (East South Central)
$ws.Range("C30").Value = 1.0900000000000001
$ws.Range("D30").Value = 11.8

# --- West South Central ---
$ws.Range("B34").Value = 10        # Row 34 (Louisiana)
$ws.Range("C34").Value = 1.23
$ws.Range("D34").Value = 9.8000000000000007

$ws.Range("B35").Value = 18        # Row 35 (Oklahoma)
$ws.Range("C35").Value = 1.81
$ws.Range("D35").Value = 21.9

# --- Mountain ---
$ws.Range("B37").Value = 3         # Row 37 (Colorado)
$ws.Range("C37").Value = 0.76
$ws.Range("D37").Value = 8.6999999999999993

$ws.Range("B38").Value = 26        # Row 38 (Idaho)
$ws.Range("D38").Value = 7

$ws.Range("B40").Value = 72        # Row 40 (Montana)
$ws.Range("D40").Value = 6

$ws.Range("B44").Value = 72        # Row 44 (New Mexico)
$ws.Range("D44").Value = 6

$ws.Range("B45").Value = 5         # Row 45 (Utah)
$ws.Range("B46").Value = 5         # Row 46 (Wyoming)

# --- Pacific Contiguous ---
$ws.Range("B50").Value = 0         # Row 50 (Pacific Contiguous)
$ws.Range("C50").Value = "--"
$ws.Range("D50").Value = "--"

$ws.Range("B57").Value = 0         # Row 57 (Washington)
$ws.Range("C57").Value = "--"
$ws.Range("D57").Value = "--"

# --- Pacific Noncontiguous ---
$ws.Range("B59").Value = 44        # Row 59 (Pacific Noncontiguous)
$ws.Range("C59").Value = 0.48
$ws.Range("D59").Value = 10

$ws.Range("B60").Value = 44        # Row 60 (Alaska)
$ws.Range("C60").Value = 0.48
$ws.Range("D60").Value = 10

# --- U.S. Total ---
$ws.Range("B66").Value = 327       # Row 66 (U.S. Total)
$ws.Range("D66").Value = 8.6999999999999993
$ws.Range("E66").Value = 124
$ws.Range("F66").Value = 0.45
$ws.Range("G66").Value = 5.5
